# Refresh the cryptos price/volume table (values from the latest coinranking.com pull).
# Column D (Price) cells are forced to text with a leading apostrophe so Excel
# keeps the exact "1.234.56"-style formatted string instead of auto-converting
# it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '51.584.07'
$ws.Range("E2").Value = '  -0.87%  '

$ws.Range("D3").Value = "'" + '2.935.86'
$ws.Range("E3").Value = '  -2.42%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = "'" + '374.70'
$ws.Range("E5").Value = '  +5.90%  '

$ws.Range("D6").Value = "'" + '104.00'
$ws.Range("E6").Value = '  -2.58%  '

$ws.Range("D7").Value = "'" + '0.543'
$ws.Range("E7").Value = '  -2.90%  '

$ws.Range("D8").Value = "'" + '0.999'
$ws.Range("E8").Value = '  -0.24%  '

$ws.Range("D9").Value = "'" + '0.587'
$ws.Range("E9").Value = '  -4.13%  '

$ws.Range("D10").Value = "'" + '37.01'
$ws.Range("E10").Value = '  -2.70%  '

$ws.Range("E12").Value = '  -2.44%  '

$ws.Range("D13").Value = "'" + '18.37'
$ws.Range("E13").Value = '  -3.32%  '

$ws.Range("D14").Value = "'" + '3.392.22'
$ws.Range("E14").Value = '  -2.44%  '

$ws.Range("D15").Value = "'" + '7.37'
$ws.Range("E15").Value = '  -3.05%  '

$ws.Range("D16").Value = "'" + '2.929.10'
$ws.Range("E16").Value = '  -2.36%  '

$ws.Range("D17").Value = "'" + '0.931'
$ws.Range("E17").Value = '  -8.31%  '

$ws.Range("D18").Value = "'" + '51.481.99'
$ws.Range("E18").Value = '  -1.24%  '

$ws.Range("E19").Value = '  +0.24%  '

$ws.Range("D20").Value = "'" + '7.34'
$ws.Range("E20").Value = '  -1.69%  '

$ws.Range("D21").Value = "'" + '12.98'
$ws.Range("E21").Value = '  -4.17%  '

$ws.Range("D22").Value = "'" + '0.0₃0947'
$ws.Range("E22").Value = '  -2.50%  '

$ws.Range("D23").Value = "'" + '68.35'
$ws.Range("E23").Value = '  -1.09%  '

$ws.Range("D24").Value = "'" + '262.29'
$ws.Range("E24").Value = '  -0.57%  '

$ws.Range("D25").Value = "'" + '2.76'
$ws.Range("E25").Value = '  +1.14%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = "'" + '0.168'
$ws.Range("E26").Value = '  -5.09%  '

$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = "'" + '4.13'
$ws.Range("E27").Value = '  -4.73%  '

$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = "'" + '7.36'
$ws.Range("E29").Value = '  -1.80%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = "'" + '25.80'
$ws.Range("E30").Value = '  -4.57%  '

$ws.Range("E31").Value = '  +7.26%  '

$ws.Range("E32").Value = '  -5.22%  '

$ws.Range("D33").Value = "'" + '9.83'
$ws.Range("E33").Value = '  -3.55%  '

$ws.Range("E34").Value = '  -3.20%  '

$ws.Range("D35").Value = "'" + '51.12'
$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("D36").Value = "'" + '34.18'
$ws.Range("E36").Value = '  -5.11%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'" + '0.0428'
$ws.Range("E37").Value = '  -2.16%  '

$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").Value = "'" + '1.00'
$ws.Range("E38").Value = '  +0.37%  '

$ws.Range("D39").Value = "'" + '3.01'
$ws.Range("E39").Value = '  -9.22%  '

$ws.Range("D40").Value = "'" + '17.02'
$ws.Range("E40").Value = '  -3.18%  '

$ws.Range("D41").Value = "'" + '2.62'
$ws.Range("E41").Value = '  -7.51%  '

$ws.Range("E42").Value = '  -6.76%  '

$ws.Range("E43").Value = '  -2.20%  '

$ws.Range("D44").Value = "'" + '122.13'
$ws.Range("E44").Value = '  -1.68%  '

$ws.Range("D45").Value = "'" + '21.85'
$ws.Range("E45").Value = '  -6.85%  '

$ws.Range("D46").Value = "'" + '2.06'
$ws.Range("E46").Value = '  -5.94%  '

$ws.Range("D47").Value = "'" + '0.274'
$ws.Range("E47").Value = '  +12.08%  '

$ws.Range("D48").Value = "'" + '2.024.36'
$ws.Range("E48").Value = '  -4.55%  '

$ws.Range("D50").Value = "'" + '3.18'
$ws.Range("E50").Value = '  -4.78%  '

$ws.Range("D51").Value = "'" + '3.212.02'
$ws.Range("E51").Value = '  -2.76%  '
